# Automatic update of files.
# Bumps the "Förändrad" (changed) date in column C by one day for every
# data row of the "Avverkningsanmälningar" sheet (rows 2-77), matching
# the daily automatic refresh of this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$firstRow = 2
$lastRow = 77

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}
